$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update description text for existing migration measures (rows 6 and 7)
$ws.Range("D6").Value = "Nettomigration till/från ett geografisk område (arbete)."
$ws.Range("D7").Value = "Nettomigration till/från ett geografisk område (utbildning)."

# Add two new rows for per-capita migration measures
$ws.Range("A8").Value = "flytt_arbete_p_cap"
$ws.Range("A9").Value = "flytt_utbildning_p_cap"

$ws.Range("B8").Value = "Flytt arbete"
$ws.Range("B9").Value = "Flytt utbildning"

$ws.Range("C8").Value = "Flytt arbete"
$ws.Range("C9").Value = "Flytt utbildning"

$ws.Range("D8").Value = "Nettomigration till/från ett geografisk område (arbete) per capita."
$ws.Range("D9").Value = "Nettomigration till/från ett geografisk område (utbildning) per capita."

$ws.Range("D11").Select()
